$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet handles
# ---------------------------------------------------------------------------
$wsStudents    = $wb.Worksheets.Item("student_data")
$wsAdvisors    = $wb.Worksheets.Item("advisors_data")
$wsInstructors = $wb.Worksheets.Item("instructors_data")
$wsStaff       = $wb.Worksheets.Item("staff_data")
$wsCourses     = $wb.Worksheets.Item("course_data")
$wsTaken       = $wb.Worksheets.Item("taken_data")
$wsMajors      = $wb.Worksheets.Item("major_data")
$wsDept        = $wb.Worksheets.Item("departmentid")
$wsLogs        = $wb.Worksheets.Item("system_logs")

# ---------------------------------------------------------------------------
# 1) instructors_data: new instructors instructor4 / instructor5 / instructor6
# ---------------------------------------------------------------------------
$wsInstructors.Cells.Item(5, 1).Value = 4
$wsInstructors.Cells.Item(5, 2).Value = "instructor4"
$wsInstructors.Cells.Item(5, 3).Value = "pass4"
$wsInstructors.Cells.Item(5, 4).Value = 2

$wsInstructors.Cells.Item(6, 1).Value = 5
$wsInstructors.Cells.Item(6, 2).Value = "instructor5"
$wsInstructors.Cells.Item(6, 3).Value = "pass5"
$wsInstructors.Cells.Item(6, 4).Value = 3

$wsInstructors.Cells.Item(7, 1).Value = 6
$wsInstructors.Cells.Item(7, 2).Value = "instructor6"
$wsInstructors.Cells.Item(7, 3).Value = "pass6"
$wsInstructors.Cells.Item(7, 4).Value = 3

# ---------------------------------------------------------------------------
# 2) course_data: fix up the previously-blank row 9, add new course rows
#    9-15 (business courses), update instructor counts on rows 7 & 8.
# ---------------------------------------------------------------------------
$wsCourses.Cells.Item(9, 1).Value  = 8
$wsCourses.Cells.Item(9, 2).Value  = "BUL"
$wsCourses.Cells.Item(9, 3).Value  = 3320
$wsCourses.Cells.Item(9, 4).Value  = "Law And Business I"
$wsCourses.Cells.Item(9, 5).Value  = 3
$wsCourses.Cells.Item(9, 6).Value  = "Fall"
$wsCourses.Cells.Item(9, 7).Value  = 2024
$wsCourses.Cells.Item(9, 8).Value  = "14:00-15:15"
$wsCourses.Cells.Item(9, 9).Value  = "TR"
$wsCourses.Cells.Item(9, 10).Value = 3
$wsCourses.Cells.Item(9, 11).Value = 2

$wsCourses.Cells.Item(10, 1).Value  = 9
$wsCourses.Cells.Item(10, 2).Value  = "CAI"
$wsCourses.Cells.Item(10, 3).Value  = 3801
$wsCourses.Cells.Item(10, 4).Value  = "Artificial Intelligence and Business Analytics for Organizations"
$wsCourses.Cells.Item(10, 5).Value  = 3
$wsCourses.Cells.Item(10, 6).Value  = "Fall"
$wsCourses.Cells.Item(10, 7).Value  = 2024
$wsCourses.Cells.Item(10, 8).Value  = "9:30-10:45"
$wsCourses.Cells.Item(10, 9).Value  = "TR"
$wsCourses.Cells.Item(10, 10).Value = 3
$wsCourses.Cells.Item(10, 11).Value = 2

$wsCourses.Cells.Item(11, 1).Value  = 10
$wsCourses.Cells.Item(11, 2).Value  = "FIN"
$wsCourses.Cells.Item(11, 3).Value  = 3403
$wsCourses.Cells.Item(11, 4).Value  = "Principles of Finance"
$wsCourses.Cells.Item(11, 5).Value  = 3
$wsCourses.Cells.Item(11, 6).Value  = "Fall"
$wsCourses.Cells.Item(11, 7).Value  = 2024
$wsCourses.Cells.Item(11, 8).Value  = "17:00-18:15"
$wsCourses.Cells.Item(11, 9).Value  = "MW"
$wsCourses.Cells.Item(11, 10).Value = 3
$wsCourses.Cells.Item(11, 11).Value = 2

$wsCourses.Cells.Item(12, 1).Value  = 11
$wsCourses.Cells.Item(12, 2).Value  = "MAN"
$wsCourses.Cells.Item(12, 3).Value  = 3025
$wsCourses.Cells.Item(12, 4).Value  = "Principles of Management"
$wsCourses.Cells.Item(12, 5).Value  = 3
$wsCourses.Cells.Item(12, 6).Value  = "Fall"
$wsCourses.Cells.Item(12, 7).Value  = 2024
$wsCourses.Cells.Item(12, 8).Value  = "8:00-9:15"
$wsCourses.Cells.Item(12, 9).Value  = "MW"
$wsCourses.Cells.Item(12, 10).Value = 4
$wsCourses.Cells.Item(12, 11).Value = 2

$wsCourses.Cells.Item(13, 1).Value  = 12
$wsCourses.Cells.Item(13, 2).Value  = "QMB"
$wsCourses.Cells.Item(13, 3).Value  = 3302
$wsCourses.Cells.Item(13, 4).Value  = "Data Analytics for Business"
$wsCourses.Cells.Item(13, 5).Value  = 3
$wsCourses.Cells.Item(13, 6).Value  = "Fall"
$wsCourses.Cells.Item(13, 7).Value  = 2024
$wsCourses.Cells.Item(13, 8).Value  = "14:00-15:15"
$wsCourses.Cells.Item(13, 9).Value  = "MW"
$wsCourses.Cells.Item(13, 10).Value = 4
$wsCourses.Cells.Item(13, 11).Value = 2

$wsCourses.Cells.Item(14, 1).Value  = 13
$wsCourses.Cells.Item(14, 2).Value  = "RMI"
$wsCourses.Cells.Item(14, 3).Value  = 3004
$wsCourses.Cells.Item(14, 4).Value  = "Principles of Risk Management"
$wsCourses.Cells.Item(14, 5).Value  = 3
$wsCourses.Cells.Item(14, 6).Value  = "Fall"
$wsCourses.Cells.Item(14, 7).Value  = 2024
$wsCourses.Cells.Item(14, 8).Value  = "9:00-10:50"
$wsCourses.Cells.Item(14, 9).Value  = "F"
$wsCourses.Cells.Item(14, 10).Value = 4
$wsCourses.Cells.Item(14, 11).Value = 2

$wsCourses.Cells.Item(15, 1).Value  = 14
$wsCourses.Cells.Item(15, 2).Value  = "MAN"
$wsCourses.Cells.Item(15, 3).Value  = 4504
$wsCourses.Cells.Item(15, 4).Value  = "Operations and Supply Chain Management"
$wsCourses.Cells.Item(15, 5).Value  = 3
$wsCourses.Cells.Item(15, 6).Value  = "Fall"
$wsCourses.Cells.Item(15, 7).Value  = 2024
$wsCourses.Cells.Item(15, 8).Value  = "15:30-16:45"
$wsCourses.Cells.Item(15, 9).Value  = "TR"
$wsCourses.Cells.Item(15, 10).Value = 3
$wsCourses.Cells.Item(15, 11).Value = 2

$wsCourses.Cells.Item(7, 10).Value = 2
$wsCourses.Cells.Item(8, 10).Value = 2

# ---------------------------------------------------------------------------
# 3) departmentid: new Arts / SOC department
# ---------------------------------------------------------------------------
$wsDept.Cells.Item(4, 1).Value = 3
$wsDept.Cells.Item(4, 2).Value = "Arts"
$wsDept.Cells.Item(4, 3).Value = "SOC"
$wsDept.Cells.Item(4, 4).Value = 330

# ---------------------------------------------------------------------------
# 4) major_data: Geology / Photography / Astronomy
# ---------------------------------------------------------------------------
$wsMajors.Cells.Item(8, 1).Value  = 7
$wsMajors.Cells.Item(8, 2).Value  = "Geology"
$wsMajors.Cells.Item(8, 3).Value  = 3

$wsMajors.Cells.Item(9, 1).Value  = 8
$wsMajors.Cells.Item(9, 2).Value  = "Photography"
$wsMajors.Cells.Item(9, 3).Value  = 3

$wsMajors.Cells.Item(10, 1).Value = 9
$wsMajors.Cells.Item(10, 2).Value = "Astronomy"
$wsMajors.Cells.Item(10, 3).Value = 3

# ---------------------------------------------------------------------------
# 5) student_data: new students student11..student17 (A/D columns, then all
#    of column B top-to-bottom, then C18, then C12..C17 -- matches the order
#    the data was actually keyed in).
# ---------------------------------------------------------------------------
$wsStudents.Cells.Item(12, 1).Value = 11
$wsStudents.Cells.Item(12, 4).Value = 7
$wsStudents.Cells.Item(13, 1).Value = 12
$wsStudents.Cells.Item(13, 4).Value = 7
$wsStudents.Cells.Item(14, 1).Value = 13
$wsStudents.Cells.Item(14, 4).Value = 7
$wsStudents.Cells.Item(15, 1).Value = 14
$wsStudents.Cells.Item(15, 4).Value = 8
$wsStudents.Cells.Item(16, 1).Value = 15
$wsStudents.Cells.Item(16, 4).Value = 8
$wsStudents.Cells.Item(17, 1).Value = 16
$wsStudents.Cells.Item(17, 4).Value = 9
$wsStudents.Cells.Item(18, 1).Value = 17
$wsStudents.Cells.Item(18, 4).Value = 9

$wsStudents.Cells.Item(12, 2).Value = "student11"
$wsStudents.Cells.Item(13, 2).Value = "student12"
$wsStudents.Cells.Item(14, 2).Value = "student13"
$wsStudents.Cells.Item(15, 2).Value = "student14"
$wsStudents.Cells.Item(16, 2).Value = "student15"
$wsStudents.Cells.Item(17, 2).Value = "student16"
$wsStudents.Cells.Item(18, 2).Value = "student17"

$wsStudents.Cells.Item(18, 3).Value = "pass17"
$wsStudents.Cells.Item(12, 3).Value = "pass11"
$wsStudents.Cells.Item(13, 3).Value = "pass12"
$wsStudents.Cells.Item(14, 3).Value = "pass13"
$wsStudents.Cells.Item(15, 3).Value = "pass14"
$wsStudents.Cells.Item(16, 3).Value = "pass15"
$wsStudents.Cells.Item(17, 3).Value = "pass16"

# ---------------------------------------------------------------------------
# 6) advisors_data: advisor4 / advisor5
# ---------------------------------------------------------------------------
$wsAdvisors.Cells.Item(5, 1).Value = 4
$wsAdvisors.Cells.Item(5, 2).Value = "advisor4"
$wsAdvisors.Cells.Item(5, 3).Value = "pass4"
$wsAdvisors.Cells.Item(5, 4).Value = 3

$wsAdvisors.Cells.Item(6, 1).Value = 5
$wsAdvisors.Cells.Item(6, 2).Value = "advisor5"
$wsAdvisors.Cells.Item(6, 3).Value = "pass5"
$wsAdvisors.Cells.Item(6, 4).Value = 3

# ---------------------------------------------------------------------------
# 7) staff_data: staff4 / staff5
# ---------------------------------------------------------------------------
$wsStaff.Cells.Item(5, 1).Value = 4
$wsStaff.Cells.Item(5, 2).Value = "staff4"
$wsStaff.Cells.Item(5, 3).Value = "pass4"
$wsStaff.Cells.Item(5, 4).Value = 3

$wsStaff.Cells.Item(6, 1).Value = 5
$wsStaff.Cells.Item(6, 2).Value = "staff5"
$wsStaff.Cells.Item(6, 3).Value = "pass5"
$wsStaff.Cells.Item(6, 4).Value = 3

# ---------------------------------------------------------------------------
# 8) course_data: more new course rows 16-22 (geology/earth-science/
#    astronomy/photography courses), added after the student/advisor/staff
#    tables above.
# ---------------------------------------------------------------------------
$wsCourses.Cells.Item(16, 1).Value  = 15
$wsCourses.Cells.Item(16, 2).Value  = "GLY"
$wsCourses.Cells.Item(16, 3).Value  = 3866
$wsCourses.Cells.Item(16, 4).Value  = "Computational Geology"
$wsCourses.Cells.Item(16, 5).Value  = 3
$wsCourses.Cells.Item(16, 6).Value  = "Fall"
$wsCourses.Cells.Item(16, 7).Value  = 2024
$wsCourses.Cells.Item(16, 8).Value  = "14:00-15:15"
$wsCourses.Cells.Item(16, 9).Value  = "TR"
$wsCourses.Cells.Item(16, 10).Value = 3
$wsCourses.Cells.Item(16, 11).Value = 3

$wsCourses.Cells.Item(17, 1).Value  = 16
$wsCourses.Cells.Item(17, 2).Value  = "GLY"
$wsCourses.Cells.Item(17, 3).Value  = 2090
$wsCourses.Cells.Item(17, 4).Value  = "History of Life"
$wsCourses.Cells.Item(17, 5).Value  = 3
$wsCourses.Cells.Item(17, 6).Value  = "Fall"
$wsCourses.Cells.Item(17, 7).Value  = 2024
$wsCourses.Cells.Item(17, 8).Value  = "9:30-10:45"
$wsCourses.Cells.Item(17, 9).Value  = "TR"
$wsCourses.Cells.Item(17, 10).Value = 3
$wsCourses.Cells.Item(17, 11).Value = 3

$wsCourses.Cells.Item(18, 1).Value  = 17
$wsCourses.Cells.Item(18, 2).Value  = "ESC"
$wsCourses.Cells.Item(18, 3).Value  = 2000
$wsCourses.Cells.Item(18, 4).Value  = "Introduction Earth Science"
$wsCourses.Cells.Item(18, 5).Value  = 3
$wsCourses.Cells.Item(18, 6).Value  = "Fall"
$wsCourses.Cells.Item(18, 7).Value  = 2024
$wsCourses.Cells.Item(18, 8).Value  = "17:00-18:15"
$wsCourses.Cells.Item(18, 9).Value  = "MW"
$wsCourses.Cells.Item(18, 10).Value = 3
$wsCourses.Cells.Item(18, 11).Value = 3

$wsCourses.Cells.Item(19, 1).Value  = 18
$wsCourses.Cells.Item(19, 2).Value  = "AST"
$wsCourses.Cells.Item(19, 3).Value  = 2002
$wsCourses.Cells.Item(19, 4).Value  = "Descriptive Astronomy"
$wsCourses.Cells.Item(19, 5).Value  = 3
$wsCourses.Cells.Item(19, 6).Value  = "Fall"
$wsCourses.Cells.Item(19, 7).Value  = 2024
$wsCourses.Cells.Item(19, 8).Value  = "8:00-9:15"
$wsCourses.Cells.Item(19, 9).Value  = "MW"
$wsCourses.Cells.Item(19, 10).Value = 4
$wsCourses.Cells.Item(19, 11).Value = 3

$wsCourses.Cells.Item(20, 1).Value  = 19
$wsCourses.Cells.Item(20, 2).Value  = "AST"
$wsCourses.Cells.Item(20, 3).Value  = 2004
$wsCourses.Cells.Item(20, 4).Value  = "Stellar Astronomy and Cosmology"
$wsCourses.Cells.Item(20, 5).Value  = 3
$wsCourses.Cells.Item(20, 6).Value  = "Fall"
$wsCourses.Cells.Item(20, 7).Value  = 2024
$wsCourses.Cells.Item(20, 8).Value  = "14:00-15:15"
$wsCourses.Cells.Item(20, 9).Value  = "MW"
$wsCourses.Cells.Item(20, 10).Value = 4
$wsCourses.Cells.Item(20, 11).Value = 3

$wsCourses.Cells.Item(21, 1).Value  = 20
$wsCourses.Cells.Item(21, 2).Value  = "PGY"
$wsCourses.Cells.Item(21, 3).Value  = "2401C"
$wsCourses.Cells.Item(21, 4).Value  = "Beginning Photography"
$wsCourses.Cells.Item(21, 5).Value  = 3
$wsCourses.Cells.Item(21, 6).Value  = "Fall"
$wsCourses.Cells.Item(21, 7).Value  = 2024
$wsCourses.Cells.Item(21, 8).Value  = "9:00-10:50"
$wsCourses.Cells.Item(21, 9).Value  = "F"
$wsCourses.Cells.Item(21, 10).Value = 4
$wsCourses.Cells.Item(21, 11).Value = 3

$wsCourses.Cells.Item(22, 1).Value  = 21
$wsCourses.Cells.Item(22, 2).Value  = "PGY"
$wsCourses.Cells.Item(22, 3).Value  = "4420C"
$wsCourses.Cells.Item(22, 4).Value  = "Advanced Photography"
$wsCourses.Cells.Item(22, 5).Value  = 3
$wsCourses.Cells.Item(22, 6).Value  = "Fall"
$wsCourses.Cells.Item(22, 7).Value  = 2024
$wsCourses.Cells.Item(22, 8).Value  = "15:30-16:45"
$wsCourses.Cells.Item(22, 9).Value  = "TR"
$wsCourses.Cells.Item(22, 10).Value = 3
$wsCourses.Cells.Item(22, 11).Value = 3

# ---------------------------------------------------------------------------
# 9) taken_data: new grade rows for the newly added students/courses
# ---------------------------------------------------------------------------
$takenRows = @(
    @(5, 9, 2.7),
    @(5, 10, 2.2999999999999998),
    @(6, 11, 3),
    @(6, 12, 4),
    @(6, 13, 3.3),
    @(7, 14, 2),
    @(7, 13, 3),
    @(7, 12, 4),
    @(8, 1, 3.3),
    @(8, 2, 4),
    @(9, 3, 4),
    @(9, 4, 3),
    @(9, 5, 3.3),
    @(10, 10, 2.7),
    @(11, 15, 2),
    @(12, 16, 2),
    @(13, 17, 1),
    @(14, 18, 3),
    @(15, 19, 3.3),
    @(16, 20, 4),
    @(17, 21, 3.6)
)

$r = 15
foreach ($rowData in $takenRows) {
    $wsTaken.Cells.Item($r, 1).Value = $rowData[0]
    $wsTaken.Cells.Item($r, 2).Value = $rowData[1]
    $wsTaken.Cells.Item($r, 3).Value = $rowData[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 10) sheet view / selection bookkeeping to mirror the final saved state
# ---------------------------------------------------------------------------
$wsStudents.Range("F11").Select()
$wsAdvisors.Range("F17").Select()
$wsInstructors.Range("G9").Select()
$wsStaff.Range("H8").Select()
$wsCourses.Range("F26").Select()
$wsTaken.Range("H31").Select()
$wsMajors.Range("I18").Select()
$wsDept.Range("I24").Select()
$wsLogs.Range("F13").Select()

$wsMajors.Select()
$wsStudents.Cells.Item(1,1).Select()
